# The deck's design was switched from the custom "Integral" theme to the
# default "Office Theme" (the standard Office color palette). In the saved
# OOXML this shows up as the in-use theme part (ppt/theme/theme1.xml, the
# one wired to the slide master) taking on the Office Theme's 12 colour
# slots (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) - its font scheme and
# format scheme were already identical between the two built-in themes, so
# only the colour scheme actually changes.
#
# Apply this by rewriting the active theme's colour scheme through the
# slide's ThemeColorScheme collection (items are ordered dk1, lt1, dk2,
# lt2, accent1..accent6, hlink, folHlink and are shared by every slide
# since they all point at the one slide master's theme, so a single pass
# on slide 1 retints the whole deck).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$themeColors = $s.ThemeColorScheme

# Target "Office Theme" palette, expressed as COM long RGB values
# (0xBBGGRR, i.e. red in the low byte) for each of the 12 theme colour
# slots, in clrScheme order.
$officeThemeRGB = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeRGB[$i - 1]
}
